# إضافة حدث جديد في Card9 — append a new service-event row (row 27) to the
# Card9 sheet, and backfill row 26's empty "nan" placeholder cells
# (B:K, M) to match the sheet's established "nan" convention for blanks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# --- Row 26: the previously blank columns get the sheet's "nan" placeholder ---
foreach ($col in 2..11) {
    $ws.Cells.Item(26, $col).Value = "nan"
}
$ws.Cells.Item(26, 13).Value = "nan"

# --- Row 27: brand-new service event ---
$ws.Cells.Item(27, 1).Value = "9"
$ws.Cells.Item(27, 12).Value = "19\7\2025"
$ws.Cells.Item(27, 14).Value = "تم تغير الجرائد الخلفيه (1_5_8)"
$ws.Cells.Item(27, 15).Value = "الخبير"
